# Auto-generated: apply scraped market-data value updates across the workbook
# (scheduled runner refresh of currentAveragePrice / LevePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1377.75
$ws.Range("I28").Value = 600.5333000000001
$ws.Range("K28").Value = 600.5333000000001
$ws.Range("M28").Value = -115.5333000000001
$ws.Range("H34").Value = 13431
$ws.Range("I34").Value = 2603.75
$ws.Range("K34").Value = 2603.75
$ws.Range("M34").Value = -2400.75
$ws.Range("H36").Value = 13431
$ws.Range("I36").Value = 2603.75
$ws.Range("K36").Value = 2603.75
$ws.Range("M36").Value = -1888.75
$ws.Range("H70").Value = 1026.2222
$ws.Range("I70").Value = 773
$ws.Range("J70").Value = 1532.6666
$ws.Range("K70").Value = 2319
$ws.Range("L70").Value = 4597.9998
$ws.Range("M70").Value = -2049
$ws.Range("N70").Value = -5137.9998
$ws.Range("H73").Value = 1026.2222
$ws.Range("I73").Value = 773
$ws.Range("J73").Value = 1532.6666
$ws.Range("K73").Value = 2319
$ws.Range("L73").Value = 4597.9998
$ws.Range("M73").Value = -1383
$ws.Range("N73").Value = -6469.9998
$ws.Range("H112").Value = 7321
$ws.Range("J112").Value = 7777.1577
$ws.Range("L112").Value = 23331.4731
$ws.Range("N112").Value = -25547.4731
$ws.Range("H116").Value = 5441.125
$ws.Range("I116").Value = 4463
$ws.Range("J116").Value = 5666.846
$ws.Range("K116").Value = 4463
$ws.Range("L116").Value = 5666.846
$ws.Range("M116").Value = -1021
$ws.Range("N116").Value = -12550.846
$ws.Range("H132").Value = 71435730
$ws.Range("I132").Value = 90917570
$ws.Range("J132").Value = 2324.3333
$ws.Range("K132").Value = 272752710
$ws.Range("L132").Value = 6972.999899999999
$ws.Range("M132").Value = -272750180
$ws.Range("N132").Value = -12032.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11157.593
$ws.Range("I134").Value = 11463.15
$ws.Range("J134").Value = 10284.571
$ws.Range("K134").Value = 34389.45
$ws.Range("L134").Value = 30853.713
$ws.Range("M134").Value = -31854.45
$ws.Range("N134").Value = -35923.713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7922
$ws.Range("I86").Value = 6336.8335
$ws.Range("J86").Value = 10035.556
$ws.Range("K86").Value = 6336.8335
$ws.Range("L86").Value = 10035.556
$ws.Range("M86").Value = -5213.8335
$ws.Range("N86").Value = -12281.556
$ws.Range("H89").Value = 7922
$ws.Range("I89").Value = 6336.8335
$ws.Range("J89").Value = 10035.556
$ws.Range("K89").Value = 31684.1675
$ws.Range("L89").Value = 50177.78
$ws.Range("M89").Value = -26068.1675
$ws.Range("N89").Value = -61409.78
$ws.Range("H99").Value = 5731.9165
$ws.Range("I99").Value = 4635.6665
$ws.Range("J99").Value = 6828.1665
$ws.Range("K99").Value = 4635.6665
$ws.Range("L99").Value = 6828.1665
$ws.Range("M99").Value = -3137.6665
$ws.Range("N99").Value = -9824.166499999999
$ws.Range("H106").Value = 50671
$ws.Range("J106").Value = 50671
$ws.Range("L106").Value = 50671
$ws.Range("N106").Value = -53195
$ws.Range("H126").Value = 5731.9165
$ws.Range("I126").Value = 4635.6665
$ws.Range("J126").Value = 6828.1665
$ws.Range("K126").Value = 13906.9995
$ws.Range("L126").Value = 20484.4995
$ws.Range("M126").Value = -11436.9995
$ws.Range("N126").Value = -25424.4995
$ws.Range("H134").Value = 10545.37
$ws.Range("I134").Value = 8785.4
$ws.Range("J134").Value = 12745.333
$ws.Range("K134").Value = 26356.2
$ws.Range("L134").Value = 38235.999
$ws.Range("M134").Value = -23821.2
$ws.Range("N134").Value = -43305.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 63622
$ws.Range("J12").Value = 129.5
$ws.Range("L12").Value = 388.5
$ws.Range("N12").Value = -734.5
$ws.Range("H25").Value = 499.35294
$ws.Range("I25").Value = 405.5625
$ws.Range("K25").Value = 1216.6875
$ws.Range("M25").Value = -1047.6875
$ws.Range("H28").Value = 2680.8333
$ws.Range("I28").Value = 2542.5
$ws.Range("K28").Value = 7627.5
$ws.Range("M28").Value = -7395.5
$ws.Range("H30").Value = 499.35294
$ws.Range("I30").Value = 405.5625
$ws.Range("K30").Value = 1216.6875
$ws.Range("M30").Value = -1114.6875
$ws.Range("H92").Value = 2272.1428
$ws.Range("J92").Value = 1640.6
$ws.Range("L92").Value = 4921.799999999999
$ws.Range("N92").Value = -7417.799999999999
$ws.Range("H97").Value = 989.3333
$ws.Range("J97").Value = 989.3333
$ws.Range("L97").Value = 2967.9999
$ws.Range("N97").Value = -3959.9999
$ws.Range("H132").Value = 1679.6666
$ws.Range("I132").Value = 1050.2858
$ws.Range("J132").Value = 2230.375
$ws.Range("K132").Value = 9452.572200000001
$ws.Range("L132").Value = 20073.375
$ws.Range("M132").Value = -6922.572200000001
$ws.Range("N132").Value = -25133.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1640.4
$ws.Range("I3").Value = 1925.5
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 1925.5
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = -1809.5
$ws.Range("N3").Value = -732
$ws.Range("H113").Value = 15154989
$ws.Range("I113").Value = 41667470
$ws.Range("K113").Value = 41667470
$ws.Range("M113").Value = -41665300
$ws.Range("H122").Value = 408655.28
$ws.Range("I122").Value = 689099.7
$ws.Range("J122").Value = 3568.889
$ws.Range("K122").Value = 2067299.1
$ws.Range("L122").Value = 10706.667
$ws.Range("M122").Value = -2064849.1
$ws.Range("N122").Value = -15606.667
$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -64900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2536
$ws.Range("I68").Value = 3222
$ws.Range("J68").Value = 1850
$ws.Range("K68").Value = 3222
$ws.Range("L68").Value = 1850
$ws.Range("M68").Value = -2473
$ws.Range("N68").Value = -3348
$ws.Range("H71").Value = 2536
$ws.Range("I71").Value = 3222
$ws.Range("J71").Value = 1850
$ws.Range("K71").Value = 16110
$ws.Range("L71").Value = 9250
$ws.Range("M71").Value = -12366
$ws.Range("N71").Value = -16738
$ws.Range("H100").Value = 4333.3335
$ws.Range("I100").Value = 4000
$ws.Range("J100").Value = 7000
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = -3459
$ws.Range("N100").Value = -8082
$ws.Range("H109").Value = 59994
$ws.Range("J109").Value = 59994
$ws.Range("L109").Value = 59994
$ws.Range("N109").Value = -62768
$ws.Range("H132").Value = 15951.6875
$ws.Range("I132").Value = 18094.46
$ws.Range("J132").Value = 6666.3335
$ws.Range("K132").Value = 54283.38
$ws.Range("L132").Value = 19999.0005
$ws.Range("M132").Value = -51753.38
$ws.Range("N132").Value = -25059.0005
$ws.Range("H136").Value = 37683.2
$ws.Range("I136").Value = 50928.19
$ws.Range("K136").Value = 152784.57
$ws.Range("M136").Value = -150234.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 14995
$ws.Range("J25").Value = 14995
$ws.Range("L25").Value = 14995
$ws.Range("N25").Value = -15581
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H109").Value = 76661.664
$ws.Range("J109").Value = 76661.664
$ws.Range("L109").Value = 76661.664
$ws.Range("N109").Value = -79435.664
$ws.Range("H113").Value = 949.875
$ws.Range("I113").Value = 785.2308
$ws.Range("K113").Value = 2355.6924
$ws.Range("M113").Value = -185.6923999999999
$ws.Range("H122").Value = 3191.8572
$ws.Range("I122").Value = 2522.5557
$ws.Range("J122").Value = 4396.6
$ws.Range("K122").Value = 7567.6671
$ws.Range("L122").Value = 13189.8
$ws.Range("M122").Value = -5117.6671
$ws.Range("N122").Value = -18089.8

